$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the metrics for the last quarter row (row 28) with refreshed data
$ws.Range("C28").Value = 369
$ws.Range("D28").Value = 41
$ws.Range("E28").Value = 328
$ws.Range("F28").Value = 6.386292834890965
